# Updated cryptos list on Tue Jul 30 13:08:52 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.716.48"
$ws.Range("E2").Value = "  -4.57%  "
$ws.Range("D3").Value = "3.354.77"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.65"
$ws.Range("E5").Value = "  -3.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.74"
$ws.Range("E6").Value = "  -5.53%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  +2.90%  "
$ws.Range("E9").Value = "  -3.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  -1.68%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.406"
$ws.Range("E11").Value = "  -3.16%  "
$ws.Range("D12").Value = "3.937.88"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.88"
$ws.Range("E14").Value = "  -6.31%  "
$ws.Range("D15").Value = "66.823.47"
$ws.Range("E15").Value = "  -4.34%  "
$ws.Range("E16").Value = "  -2.51%  "
$ws.Range("D17").Value = "3.356.14"
$ws.Range("E17").Value = "  -1.76%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "440.75"
$ws.Range("E18").Value = "  -3.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.60"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.61"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("E25").Value = "  -5.09%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.13"
$ws.Range("E27").Value = "  -3.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("E29").Value = "  -3.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.88"
$ws.Range("E30").Value = "  -2.29%  "
$ws.Range("E31").Value = "  -4.99%  "
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.82"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "161.71"
$ws.Range("E36").Value = "  -5.55%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "28.03"
$ws.Range("E37").Value = "  +2.06%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.79"
$ws.Range("E38").Value = "  -8.08%  "
$ws.Range("D39").Value = "2.828.90"
$ws.Range("E39").Value = "  +3.48%  "
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("E41").Value = "  -3.53%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.21"
$ws.Range("E42").Value = "  -5.66%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.45"
$ws.Range("E43").Value = "  -1.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0668"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "24.34"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.35"
$ws.Range("E46").Value = "  -7.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "327.72"
$ws.Range("E47").Value = "  -3.09%  "
$ws.Range("E48").Value = "  -4.17%  "
$ws.Range("E49").Value = "  +1.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.980"
$ws.Range("E50").Value = "  -3.66%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.17"
$ws.Range("E51").Value = "  -2.40%  "
